$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9226154088973999
$ws.Range("B1").Value = 2.023155689239502
$ws.Range("C1").Value = 8.864248275756836
$ws.Range("D1").Value = 1.829176783561707
$ws.Range("E1").Value = 1.426509141921997
